# Multiple Tabs task in progress
# - Clear the (stale) "Passed" status in C4/C5
# - Reword the B9 reply
# - Leave selection on B9 (last edited cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()

$ws.Range("B9").Value = "2 Days After tomorrow"

$ws.Range("B9").Select() | Out-Null
